# Apply price/volume/coin-name updates to Sheet1 as captured in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) hold text values that look numeric
# (e.g. "68.190.89", "0.0000294", "  +1.24%  "). Force the cell format to
# Text before assigning so Excel does not silently coerce them into numbers
# or dates and strip formatting such as trailing zeros / leading spaces.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.190.89'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.24%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.633.66'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.47%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.54%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '196.16'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +5.74%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '575.70'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -1.91%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.628.00'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -5.07%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.92%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.07%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.678'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.29%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +4.97%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '55.98'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +3.32%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000294'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +16.99%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.11'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.29%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.208.47'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.52%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.631.06'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.64%  '
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.19%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.54'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +2.24%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '68.136.20'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.67%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.99%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.15%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '402.83'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +1.89%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +22.38%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.23'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -1.68%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '86.11'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.29%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +3.24%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +2.73%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +6.84%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +1.25%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.13'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +19.12%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.15'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +1.95%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '31.72'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.55%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '696.32'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +18.54%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +2.54%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +5.87%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '64.77'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -2.26%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '42.73'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +2.13%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +10.89%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.12%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +8.57%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.86'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +18.77%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.136'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +2.32%  '
$ws.Range('B43').Value = 'ThetaToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.13'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +12.34%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.168.26'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +17.91%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.19%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.90'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +24.90%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0423'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.84%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.81%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +5.59%  '
$ws.Range('B50').Value = 'ApeXProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.11'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.57%  '
$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '142.54'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.94%  '
